$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=16; C='1047444936'; D='GREISTON PIMENTEL URRUTIA'; E='2503'; F=31200; G=1300000},
    @{Row=17; C='1047444936'; D='GREISTON PIMENTEL URRUTIA'; E='2502'; F=52000; G=1300000},
    @{Row=18; C='1047444936'; D='GREISTON PIMENTEL URRUTIA'; E='2501'; F=52000; G=1300000},
    @{Row=19; C='1047444936'; D='GREISTON PIMENTEL URRUTIA'; E='2412'; F=52000; G=1300000},
    @{Row=20; C='1047444936'; D='GREISTON PIMENTEL URRUTIA'; E='2411'; F=52000; G=1300000},
    @{Row=21; C='1047444936'; D='GREISTON PIMENTEL URRUTIA'; E='2410'; F=52000; G=1300000},
    @{Row=22; C='1047444936'; D='GREISTON PIMENTEL URRUTIA'; E='2409'; F=52000; G=1300000},
    @{Row=23; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2503'; F=27840; G=1160000},
    @{Row=24; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2502'; F=46400; G=1160000},
    @{Row=25; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2501'; F=46400; G=1160000},
    @{Row=26; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2412'; F=46400; G=1160000},
    @{Row=27; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2411'; F=46400; G=1160000},
    @{Row=28; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2410'; F=46400; G=1160000},
    @{Row=29; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2409'; F=46400; G=1160000},
    @{Row=30; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2408'; F=46400; G=1160000},
    @{Row=31; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2407'; F=46400; G=1160000},
    @{Row=32; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2406'; F=46400; G=1160000},
    @{Row=33; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2405'; F=46400; G=1160000},
    @{Row=34; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2404'; F=46400; G=1160000},
    @{Row=35; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2403'; F=46400; G=1160000},
    @{Row=36; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2402'; F=46400; G=1160000},
    @{Row=37; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2401'; F=46400; G=1160000},
    @{Row=38; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2312'; F=46400; G=1160000},
    @{Row=39; C='19335385'; D='EDUARDO MEDINA ROMERO'; E='2311'; F=46400; G=1160000},
    @{Row=40; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2503'; F=31200; G=1300000},
    @{Row=41; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2502'; F=52000; G=1300000},
    @{Row=42; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2501'; F=52000; G=1300000},
    @{Row=43; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2412'; F=52000; G=1300000},
    @{Row=44; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2411'; F=52000; G=1300000},
    @{Row=45; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2410'; F=52000; G=1300000},
    @{Row=46; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2409'; F=52000; G=1300000},
    @{Row=47; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2408'; F=52000; G=1300000},
    @{Row=48; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2407'; F=52000; G=1300000},
    @{Row=49; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2406'; F=52000; G=1300000},
    @{Row=50; C='15324579'; D='HECTOR NICILAS FERNANDEZ MUNERA'; E='2405'; F=52000; G=1300000}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

Write-Host "Done updating rows 16-50"
